$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

function Replace-ParagraphWithPlainText($paraIndex, $newText) {
    # Deletes the whole paragraph (including any hyperlinks / multi-run
    # content) and inserts a brand new single-run paragraph with the
    # same (BodyText) style right before what used to follow it. Net
    # paragraph count is unchanged, so indices after this call are
    # identical to indices before it.
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.Delete()
    $following = $d.Paragraphs.Item($paraIndex)
    $following.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($paraIndex)
    $newPara.Range.Text = $newText
}

# ---------------------------------------------------------------------
# 1) Bookmark ids near the top / around the projects section collapse
#    back down as a natural side effect of editing the document; no
#    explicit bookmark surgery is required here.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 2) "Proyectos Destacados" section rewrite
# ---------------------------------------------------------------------

# 73: Lends -> GEOTRANS - Maquinarias y Transportes (keeps FirstParagraph style)
Replace-InParagraph 73 "Lends" "GEOTRANS - Maquinarias y Transportes"

# 74: old Lends description -> GEOTRANS description
Replace-InParagraph 74 "PWA para gestión de préstamos con autenticación de Google y sincronización en tiempo real con Firebase." "Diseñé y desarrollé un sitio web corporativo moderno y optimizado para una empresa con más de 20 años en el rubro del transporte y maquinaria pesada en el norte de Chile. Implementé una arquitectura JAMstack con Gatsby y React, integrando formularios sin servidor con AWS Lambda, API Gateway y SES. Enfocado en SEO, rendimiento, experiencia de usuario y soporte multiidioma (español/inglés), el proyecto está desplegado como PWA de alta performance."

# 75: feb. 2024 - Presente -> abr. 2020 - may. 2020
Replace-InParagraph 75 "feb. 2024 — Presente" "abr. 2020 — may. 2020"

# 76: hyperlink display text my-lends -> geotrans.vercel.app (rId40 / url untouched)
Replace-InParagraph 76 "github.com/Freddymhs/my-lends" "geotrans.vercel.app"
$p76 = $d.Paragraphs.Item(76)
$r76 = $p76.Range
$r76.Find.Execute("geotrans.vercel.app") | Out-Null
$r76.Style = "Hyperlink"

# 77: Geotrans Landing Page -> Lends
Replace-InParagraph 77 "Geotrans Landing Page" "Lends"

# 78: Geotrans landing description -> "..."
Replace-InParagraph 78 "Landing page responsiva para empresa de transporte Geotrans, con secciones de servicios, testimonios y formulario de contacto integrado." "..."

# 79: [github.com/Freddymhs/Hosting-GeotransWeb] hyperlink paragraph -> plain date range
Replace-ParagraphWithPlainText 79 "feb. 2024 — abr. 2024"

# 80: Mappy -> Bot Router
Replace-InParagraph 80 "Mappy" "Bot Router"

# 81: Mappy description -> "..."
Replace-InParagraph 81 "Aplicación web que extrae eventos públicos de Instagram mediante scraping y los visualiza en un mapa 3D interactivo en tiempo real usando Three.js y React." "..."

# 82: [github.com/Freddymhs/mapy-front] hyperlink paragraph -> plain date range
Replace-ParagraphWithPlainText 82 "jul. 2024 — ago. 2024"

# 83: Agendo -> Gestor APP
Replace-InParagraph 83 "Agendo" "Gestor APP"

# 84: Agendo description -> "..."
Replace-InParagraph 84 "nestjs , Sistema de gestión de citas en tiempo real para clínica de rehabilitación, con autenticación, notificaciones y accesibilidad para personas con movilidad reducida." "..."

# 85: [github.com/Freddymhs/agendo] hyperlink paragraph -> plain date range
Replace-ParagraphWithPlainText 85 "feb. 2024 — ene. 2025"

# 86-88: Geographical Info project entirely removed (delete top-down so
# indices of the not-yet-deleted paragraphs stay put).
$d.Paragraphs.Item(88).Range.Delete()
$d.Paragraphs.Item(87).Range.Delete()
$d.Paragraphs.Item(86).Range.Delete()

Write-Output "done"
